$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15 (pushes rows 15-25 down to 16-26), reusing the
# formatting that Excel carries along on a row insert.
$ws.Rows(15).Insert() | Out-Null

# The freshly inserted row 15 starts out essentially blank / default
# styled; clone the (correct, original) formatting now sitting in row 16
# (the old row 15, pushed down) so row 15 matches every other data row.
$ws.Range("A16:Q16").Copy() | Out-Null
$ws.Range("A15:Q15").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Row height for the new data row matches the other 24.75pt data rows.
$ws.Rows(15).RowHeight = 24.75

# Populate the new row: GLYCERIN INFANTILE 10 SUPP. (GLAXO), inserted
# alphabetically between FERROTRON 30 CAPS and OXITROPIL 1200 MG 60 TAB.
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "GLYCERIN INFANTILE 10 SUPP. (GLAXO)"
$ws.Range("H15").Value = "35:0"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "19.00"
$ws.Range("P15").Value = "19.0000"
$ws.Range("Q15").Value = "1:0"

# Renumber the sequence column for the rows that were pushed down
# (old #9..#17 in rows 16-24 become #10..#18).
for ($r = 16; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# Update the running total (old 946.215 + the new item's 19.0000).
$ws.Range("P25").Value = 965.215

# Refresh the generated timestamp footer text.
$ws.Range("A26").Value = "Monday, 8 September, 2025 11:45 AM"
